$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 becomes FraxShare (was Algorand), Row 43 becomes Algorand (was FraxShare)
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'9.06"
$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.202"
$ws.Range("E43").Value = "  +15.43%  "

# Update Price (D) and Volume(1h) (E) columns for remaining rows
$ws.Range("D2").Value = "43.881.15"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.353.16"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +5.66%  "
$ws.Range("D6").Value = "'242.21"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").Value = "'77.25"
$ws.Range("E7").Value = "  +4.83%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +21.60%  "
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("D11").Value = "'57.38"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'33.89"
$ws.Range("E12").Value = "  +22.11%  "
$ws.Range("D13").Value = "'7.58"
$ws.Range("E13").Value = "  +19.37%  "
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").Value = "2.703.11"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "'16.99"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "'0.927"
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "2.350.75"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "43.760.94"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").Value = "'77.62"
$ws.Range("D23").Value = "'256.47"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("D26").Value = "'11.04"
$ws.Range("E26").Value = "  +8.49%  "
$ws.Range("E27").Value = "  -5.34%  "
$ws.Range("D28").Value = "'1.79"
$ws.Range("E28").Value = "  +16.56%  "
$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").Value = "'23.08"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'174.93"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("E33").Value = "  +6.07%  "
$ws.Range("D34").Value = "'0.0758"
$ws.Range("E34").Value = "  +7.94%  "
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("D36").Value = "'5.41"
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").Value = "'3.80"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'2.43"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("E40").Value = "  +7.83%  "
$ws.Range("D41").Value = "'19.55"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +6.97%  "
$ws.Range("D46").Value = "'2.54"
$ws.Range("E46").Value = "  +12.55%  "
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("D49").Value = "'101.83"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "'55.39"
$ws.Range("E51").Value = "  +8.34%  "
